$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells we touch so numeric-looking strings
# (e.g. "0.9995", "103.00") are preserved verbatim as text, matching
# the original inlineStr cell type, instead of being coerced to numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.271.70'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.861.70'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9995'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7016'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '238.04'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08341'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +11.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3050'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.38'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08183'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.865.74'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7187'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.192'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.41'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.281.18'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.37%  '
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007913'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.18%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.794'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.43'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '237.42'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.104.08'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9997'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.471'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.15'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.014'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1454'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.14'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.993'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.434'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.97%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.431'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.484'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.065'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05217'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.15%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7076'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.001'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.660'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01853'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.719'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.150.10'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +8.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9209'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.946'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4289'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '71.01'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9993'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.00'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.780'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.001.56'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.207'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.990'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.84%  '
